$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.744.86'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.561.69'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.28'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('E6').Value = '  -2.67%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '21.89'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.782.32'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.563.69'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('E14').Value = '  -2.54%  '
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '61.46'
$ws.Range('E16').Value = '  -3.14%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.766.06'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '213.42'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.08'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.34'
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.64'
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.75'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.81'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('E30').Value = '  -4.37%  '
$ws.Range('E31').Value = '  -1.77%  '
$ws.Range('E32').Value = '  -1.92%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.382.34'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.519'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.989'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.19'
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.695.87'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.48'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0₇0981'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0493'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0947'
$ws.Range('E51').Value = '  -0.63%  '
